$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "CGTCTCCAGCATTATACAAAGCGGTTCAGGCAGAAAACCCAACAATTCAACATTTTGACGCATCAGTATTTACCGGCGAATACATTACCGGCGCGAGACG",
    "CGTCTCCGGCGATGTAGATAAAGCTTACCTTGACGCAATTGCACACGCTCGTAATGACAAAGCAAAAGCCAAAGCTGCAAAACAAGCAACCAACGAGACG",
    "CGTCTCCCCAATTTAGAAATTCATAACGAAAACTAAATTAAAGCATAAAAACAAGCCCTTGAGTCTGAGATGATTCAAGGGCTTTTTGTATTACGAGACG",
    "CGTCTCCATTAATGGCTATATAAAACGATGCCCCTTGAAACCGAAATTTCAAGGGGCAGAGCGTTGCGAGAAATTCATCTCTTATTGTAACAACGAGACG",
    "CGTCTCCACAACGTTTTTCATTTCGCAGTGCCATTACAATGATTAAATAGTATCGTTGACATTGATTGTTGTCAATATACATTTATTACAATTCGAGACG",
    "CGTCTCCAATTTATGATAAAAATTAAACGCCTTGTAAAAAAGGATTATTACAAATTATAAGAGGTGAAATTGTATGAAATTGACGCCATTAAACGAGACG",
    "CGTCTCCTAAACTATATTTTAGGTCTTGATTTAGGGATTGCTTCTGTCGGTTGGGCAGTGGTAGAGATTGATGAGCAAGAAAATCCACTAGGTCGAGACG",
    "CGTCTCCAGGTTTAATTGATGTAGGAGTACGAACATTTGATAGAGCTGAAGTGCCGAAAACAGGCGAAAGTTTGGCATTAGCTCGCCGTTTAGCGAGACG",
    "CGTCTCCTTAGCTCGTTCTGCTCGTCGTTTAGTAAAACGTCGAGCGGATCGAATTAAAAAAGCGAAGCGTTTATTAAAAGCAGAAAATATTTTCGAGACG",
    "CGTCTCCTTTTACTTTCGGCAGATGAACACTTGCCCAATGATGTTTGGCAATTACGGGTTAAAGGTTTGGATCAAAAGCTCGAACGCCAGGAACGAGACG",
    "CGTCTCCGGAATGGGCAGCGGTTTTATTGCATTTATTGAAACATCGTGGTTATTTGTCACAACGTAAAAATGAAAGCAAAAGTGAGAATAAAGCGAGACG",
    "CGTCTCCAAAGAATTAGGTGCGTTGCTTTCAGGTGTAGAAACAAACCATCAAATTTTGCAATCTGCTGAATACCGCACGCCTGCGGAAATTGCCGAGACG",
    "CGTCTCCTTGCTGTGAAGAAATTTCACGTAGAAGATGGGCATATTCGTAATCAGCGTGGTGCTTATACGCATACATTTAGCCGTTTAGATTTACGAGACG",
    "CGTCTCCTTTATTGGCGGAAATGGAATTATTGTTCCAACGCCAAACGGACTTGGGCAATCCGCACACTTCTGCAAAATTATTGGAAAATTTGACGAGACG",
    "CGTCTCCTTGACCGCTTTATTGATGTGGCAAAAGCCTGCGTTGGCGGGCGAAGCCATTTTGAAAATGCTCGGCAAATGTACCTTTGAACCCACCGAGACG",
    "CGTCTCCCCACCGAATATAAAGCGGCGAAAAATAGTTATTCGGCTGAGCGTTTTGTATGGCTGACCAAGTTGAATAATTTGCGTATTTTGGAACGAGACG",
    "CGTCTCCGGAACAAGGGGCTGAGCGTGCATTGACTGATAACGAACGTTTTGCTCTGCTCGATCAGCCTTATGAAAAAGCCAAGTTTACTTACGCGAGACG",
    "CGTCTCCTACGCCCAAGCACGTACAATGTTAGCTTTACCTGATGAAGCAATTTTTAAGGGCGTGCGTTATCAAGGCGAAGATAAAAAAGCTGTCGAGACG",
    "CGTCTCCCTGTCGAAACGAAAACTATTTTAATGGAGATGAAAGCCTATCATCAAATCCGTAAAGCATTAGAGAATGCAGATTTAAAAGCAGAACGAGACG",
    "CGTCTCCAGAATGGAATGAACTTAAAAATAATTCCGAATTGCTTGATGACATTGGCACAGCGTTTTCATTGCATAAAACTGATGAAGATATTTCGAGACG",
    "CGTCTCCATTTGCCGTTATTTAGATGGAAAATTATCGGAAAGCATATTGAATGCGTTGTTAGAAAATCTGAATTTTGACAAATTTATTCAACTCGAGACG",
    "CGTCTCCAACTTTCACTTAAAGCATTACAACAAATTTTACCGTTGATGTTGCAAGGGCAACGTTATGATGAAGCGGTTTCAGCGATTTATGGTCGAGACG",
    "CGTCTCCTGGTGATCATTATGGTAAAAAATCAGCAGAAATTAACCGCTTGTTACCAACTATTCCAGCCGATGAAATCCGCAATCCAGTAGTATCGAGACG",
    "CGTCTCCGTATTACGCACACTGACTCAAGCTCGCAAAGTGATCAATGCGGTGGTGCGATTGTATGGTTCACCTGCTCGTATTCATATTGAAACCGAGACG",
    "CGTCTCCAAACAGGACGAGAAGTGGGCAAATCTTATCAAGATCGTAAGAAACTGGAGAAACAACAGGAAGATAATCGTAAACAACGTGAAAGTCGAGACG",
    "CGTCTCCAAGTGCGGTGAAAAAATTCAAAGAATATTTTCCAAATTTCGTGGGAGAGCCAAAAGGTAAAGATATTCTAAAAATGCGTTTGTATGCGAGACG",
    "CGTCTCCTATGAGTTGCAACAAGCAAAATGTTTATATTCAGGCAAATCATTGGAATTACACCGCTTATTGGAAAAAGGCTATGTAGAAGTTGACGAGACG",
    "CGTCTCCTTGATCATGCTTTGCCGTTTTCTCGCACTTGGGATGATAGCTTTAATAATAAAGTGTTGGTGCTTGCCAATGAAAACCAAAATAAACGAGACG",
    "CGTCTCCTAAAGGCAATTTAACACCTTATGAATGGTTAGATGGCAAAAATAATAGTGAACATTGGCAAAATTTTGTCGCACGAGTACAAACTACGAGACG",
    "CGTCTCCACTAGTGGTTTCTCACATACTAAGAAACAACGTATTTTAAGTCATAAACTAGATGAAAAAGGCTTTATCGAACGTAATTTAAATGACGAGACG",
    "CGTCTCCATGATACTCGCTATGTTGCCCGTTTCTTATGTAATTTTATTGCTGACAATATGTTACTGACAGGCAAAGGCAAGCGAAAAGTGTTTCGAGACG",
    "CGTCTCCGTTTGCTTCAAATGGGCAAATTACGGCTTTATTACGTGGGCGTTGGGGTTTACAAAAAGTACGTGATGATAATGATCGCCACCACGCGAGACG",
    "CGTCTCCCACGCTTTAGATGCGGTTGTGGTTGCCTGCTCAACGGTAGTGATGCAACAGAAAATTACAAGATTTGTGAGATATGAAGAGGGTAACGAGACG",
    "CGTCTCCGTAATGTTTTCAGTGGAGAACGAATTGATCGTGAAACTGGTGAGATTATTCCTTTGCATTTCCCAAGCCCCTGGGCATTTTTTAGACGAGACG",
    "CGTCTCCTAGAGAAAATGTGGAAATTCGCATTTTTAGTGAAAATCCGAAATTAGAACTGGAAAATCGCTTACCTGATTATCCACAATATAATCCGAGACG",
    "CGTCTCCAATCACGAATTTGTTCAGCCGTTATTTGTGTCGAGAATGCCAACCCGAAAAATGACAGGGCAAGGGCATATGGAAACAGTAAAATCCGAGACG",
    "CGTCTCCAATCAGCCAAACGTTTAGATGAAGGTTTAAGTGTATTAAAAGTGCCTTTAACACAACTTAAATTGAGCGATTTAGAGCGAATGGTTCGAGACG",
    "CGTCTCCGGTTAATCGTGAGCGTGAAGTTACATTGTACGAATCTTTAAAAGCCCGTTTAGAACAATTTGGTAATGATCCAGCAAAAGCTTTTGCGAGACG",
    "CGTCTCCTTTGCCGAACCGTTCCATAAAAAAGGCGGTGCGGTGGTTAAAGCTGTGCGAGTGGAACAAACGCAAAAATCAGGCGTATTAGTGCGCGAGACG",
    "CGTCTCCTGCGTGATGGCAATGGTGTTGCGGATAATGCTTCTATGGTGCGAGTTGATGTCTTTACCAAAGGTGGCAAATACTTCCTTGTGCCACGAGACG",
    "CGTCTCCGCCAATTTACACTTGGCAAGTGGCGAAAGGGATTTTGCCAAATAAGGCGGTAACTGCTAATGTTGATGAAATTGATTGGCTTGAAACGAGACG",
    "CGTCTCCGAAATGGATGAAAGTTATCAATTTATCTTTACTATGTACCCAAATGATCTTGTCAAAGTAAAATTGAAAAAAGAAGAATTCTTTGGCGAGACG",
    "CGTCTCCTTGGTTATTACGGTGGTTTAGATCGAGCAACAGGGGCTATTGTCATAAAAGAACACGATTTAGAAAAATCCAAAGGAAAACAAGGTCGAGACG",
    "CGTCTCCAGGTATTTATCGTATTGGCGTTAAATTAGCTTTGTCATTTGAAAAATACCAAGTCGATGAACTCGGTAAAAATATCCGTCCTTGTCCGAGACG",
    "CGTCTCCTGTCGTCCAACTAAACGACAACACGTACGCTAACTGAATCCCTACACTCTTCGAGTGTGGGGATTTTTTGTATTTAAGGAAGAAATCGAGACG",
    "CGTCTCCAAATTATGACTTGGCGTAGTATTTTAATTAGCAAGGGCGGAAAACTTTCCTTGCAGAAAAATCAAATGTTGATTCAGCAAGAGGGTCGAGACG",
    "CGTCTCCGGGTAATGAATTTACTGTACCTTTGGAAGATATTGCGATTGTAGTGGTGGATAGTCGGGAAACGGTTATTACGATTCCCTTATTATCGAGACG",
    "CGTCTCCTTATCTGCTTTTGGTTTATACGGCATTACGTTTTTAACTTGTGATGAACAGTTTTTACCTTGTGGGCAATGGTTGCCATTTAATCACGAGACG",
    "CGTCTCCATCAATATCATCGACAGCTCAAAACCTTGAAATTACAGCTAGAAGCTAGCTTGCCACAAAAGAAGCAGCTTTGGCAGAAAATTGTGCGAGACG",
    "CGTCTCCTGTGCAACAGAAAATCCGAAATCAAGCGACAGTGTTGAAGATTTGCAAATTTCAAGCAGAATCCGACCGCTTGTCTAAAATGGCAGAGCAAGTAAAGAGCACGAGACG"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
